$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, even if it looks numeric,
# by temporarily forcing Text format, then restoring the original (default) style
# so we don't leave a stray number-format style behind on the cell.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range('D2').Value = '43.724.64'
$ws.Range('E2').Value = '  +1.86%  '

$ws.Range('D3').Value = '2.339.54'
$ws.Range('E3').Value = '  +1.90%  '

$ws.Range('E4').Value = '  +0.22%  '

$ws.Range('B5').Value = 'Solana'
$ws.Range('C5').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Range('D5') '109.39'
$ws.Range('E5').Value = '  +4.72%  '

$ws.Range('B6').Value = 'BNB'
$ws.Range('C6').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Range('D6') '311.53'
$ws.Range('E6').Value = '  -1.61%  '

$ws.Range('E7').Value = '  +1.06%  '

$ws.Range('E8').Value = '  -0.01%  '

Set-TextValue $ws.Range('D9') '0.619'
$ws.Range('E9').Value = '  +2.74%  '

Set-TextValue $ws.Range('D10') '41.14'
$ws.Range('E10').Value = '  +4.05%  '

Set-TextValue $ws.Range('D11') '0.0921'
$ws.Range('E11').Value = '  +1.86%  '

Set-TextValue $ws.Range('D12') '8.59'
$ws.Range('E12').Value = '  +1.33%  '

$ws.Range('E13').Value = '  -0.89%  '

Set-TextValue $ws.Range('D14') '1.01'
$ws.Range('E14').Value = '  -0.65%  '

Set-TextValue $ws.Range('D15') '15.55'
$ws.Range('E15').Value = '  +1.34%  '

$ws.Range('D16').Value = '2.685.86'
$ws.Range('E16').Value = '  +1.60%  '

$ws.Range('D17').Value = '2.323.20'
$ws.Range('E17').Value = '  +0.59%  '

$ws.Range('D18').Value = '43.601.43'
$ws.Range('E18').Value = '  +1.74%  '

Set-TextValue $ws.Range('D19') '7.60'
$ws.Range('E19').Value = '  +1.65%  '

$ws.Range('E20').Value = '  +1.56%  '

Set-TextValue $ws.Range('D21') '13.18'

Set-TextValue $ws.Range('D22') '74.33'
$ws.Range('E22').Value = '  +0.45%  '

$ws.Range('E23').Value = '  -1.43%  '

Set-TextValue $ws.Range('D24') '269.71'
$ws.Range('E24').Value = '  +2.50%  '

$ws.Range('E25').Value = '  +2.91%  '

$ws.Range('E26').Value = '  -0.09%  '

Set-TextValue $ws.Range('D27') '7.64'
$ws.Range('E27').Value = '  +7.81%  '

Set-TextValue $ws.Range('D28') '11.18'
$ws.Range('E28').Value = '  +2.31%  '

$ws.Range('E29').Value = '  +0.06%  '

Set-TextValue $ws.Range('D30') '38.80'
$ws.Range('E30').Value = '  +3.10%  '

Set-TextValue $ws.Range('D31') '22.69'
$ws.Range('E31').Value = '  +1.48%  '

Set-TextValue $ws.Range('D32') '167.93'
$ws.Range('E32').Value = '  +0.74%  '

$ws.Range('E33').Value = '  +1.54%  '

Set-TextValue $ws.Range('D34') '2.80'
$ws.Range('E34').Value = '  +8.51%  '

Set-TextValue $ws.Range('D35') '0.132'
$ws.Range('E35').Value = '  +1.00%  '

Set-TextValue $ws.Range('D36') '4.77'
$ws.Range('E36').Value = '  +4.26%  '

$ws.Range('E37').Value = '  -2.01%  '

$ws.Range('E38').Value = '  +4.30%  '

$ws.Range('E39').Value = '  -1.08%  '

Set-TextValue $ws.Range('D40') '2.86'
$ws.Range('E40').Value = '  +6.64%  '

Set-TextValue $ws.Range('D41') '1.73'
$ws.Range('E41').Value = '  +9.78%  '

Set-TextValue $ws.Range('D42') '105.08'
$ws.Range('E42').Value = '  +13.38%  '

Set-TextValue $ws.Range('D43') '71.82'
$ws.Range('E43').Value = '  +3.01%  '

Set-TextValue $ws.Range('D44') '0.237'
$ws.Range('E44').Value = '  +2.90%  '

Set-TextValue $ws.Range('D45') '13.33'
$ws.Range('E45').Value = '  +8.77%  '

$ws.Range('E46').Value = '  -0.18%  '

Set-TextValue $ws.Range('D47') '114.37'
$ws.Range('E47').Value = '  -0.07%  '

$ws.Range('D48').Value = '1.668.18'
$ws.Range('E48').Value = '  -3.23%  '

Set-TextValue $ws.Range('D49') '77.27'
$ws.Range('E49').Value = '  -3.75%  '

$ws.Range('E50').Value = '  +3.37%  '

Set-TextValue $ws.Range('D51') '8.96'
$ws.Range('E51').Value = '  +2.00%  '

